# Insert a new data row above current row 15 (Fecha 2022-08-07 / 44790),
# pushing the existing rows 15-22 down to 16-23.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Cells.Item(15, 1).Value  = 7
$ws.Cells.Item(15, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(15, 3).Value  = "Ñuble"
$ws.Cells.Item(15, 4).Value  = 44790
$ws.Cells.Item(15, 5).Value  = 16
$ws.Cells.Item(15, 6).Value  = 100112037
$ws.Cells.Item(15, 7).Value  = "Cebollín"
$ws.Cells.Item(15, 8).Value  = "Sin especificar"
$ws.Cells.Item(15, 9).Value  = "Primera"
$ws.Cells.Item(15, 10).Value = 120
$ws.Cells.Item(15, 11).Value = 8000
$ws.Cells.Item(15, 12).Value = 8500
$ws.Cells.Item(15, 13).Value = 8250
$ws.Cells.Item(15, 14).Value = '$/docena de atados'
$ws.Cells.Item(15, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(15, 16).Value = 2750
$ws.Cells.Item(15, 17).Value = 3
$ws.Cells.Item(15, 18).Value = "Hortaliza"
